# Updates the cryptos list Price (D) / Volume(1h) (E) columns per the target diff.
# D-column values that parse as plain numbers are written with a leading apostrophe
# so Excel stores them as text (preserving formatting like trailing zeros),
# matching the source data's non-numeric 'Price' strings (e.g. '26.302.76').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.302.76'
$ws.Range('E2').Value = '  +0.06%  '

$ws.Range('D3').Value = '1.689.80'
$ws.Range('E3').Value = '  +0.61%  '

$ws.Range('D4').Value = '''1.009'
$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = '''217.96'
$ws.Range('E5').Value = '  -0.19%  '

$ws.Range('D6').Value = '''0.5373'
$ws.Range('E6').Value = '  +2.42%  '

$ws.Range('D7').Value = '''1.009'
$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').Value = '''0.2730'
$ws.Range('E8').Value = '  +1.25%  '

$ws.Range('D9').Value = '''0.06440'
$ws.Range('E9').Value = '  -0.48%  '

$ws.Range('D10').Value = '''21.63'
$ws.Range('E10').Value = '  -1.67%  '

$ws.Range('D11').Value = '''0.07670'
$ws.Range('E11').Value = '  +1.83%  '

$ws.Range('D12').Value = '1.690.77'
$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('D13').Value = '''4.529'
$ws.Range('E13').Value = '  -0.08%  '

$ws.Range('D14').Value = '''0.5786'
$ws.Range('E14').Value = '  -0.31%  '

$ws.Range('D15').Value = '''0.000008368'
$ws.Range('E15').Value = '  -1.73%  '

$ws.Range('D16').Value = '''66.92'
$ws.Range('E16').Value = '  +3.40%  '

$ws.Range('D17').Value = '26.369.41'
$ws.Range('E17').Value = '  +0.16%  '

$ws.Range('D18').Value = '''4.905'
$ws.Range('E18').Value = '  -0.42%  '

$ws.Range('E20').Value = '  -0.18%  '

$ws.Range('D21').Value = '''190.09'
$ws.Range('E21').Value = '  -0.05%  '

$ws.Range('D22').Value = '''6.261'
$ws.Range('E22').Value = '  +0.83%  '

$ws.Range('E23').Value = '  +0.17%  '

$ws.Range('D24').Value = '''149.09'
$ws.Range('E24').Value = '  +2.32%  '

$ws.Range('D25').Value = '''0.1287'
$ws.Range('E25').Value = '  +3.19%  '

$ws.Range('D26').Value = '''7.859'
$ws.Range('E26').Value = '  +0.56%  '

$ws.Range('E27').Value = '  +0.33%  '

$ws.Range('D28').Value = '''0.06261'
$ws.Range('E28').Value = '  -2.98%  '

$ws.Range('D29').Value = '''1.371'

$ws.Range('D30').Value = '''1.326'
$ws.Range('E30').Value = '  -0.01%  '

$ws.Range('D31').Value = '''3.598'
$ws.Range('E31').Value = '  -0.22%  '

$ws.Range('D32').Value = '''3.582'
$ws.Range('E32').Value = '  -0.38%  '

$ws.Range('D33').Value = '''1.671'
$ws.Range('E33').Value = '  +0.26%  '

$ws.Range('D34').Value = '''1.031'
$ws.Range('E34').Value = '  +0.03%  '

$ws.Range('D35').Value = '''0.6145'
$ws.Range('E35').Value = '  -1.57%  '

$ws.Range('D36').Value = '''2.417'
$ws.Range('E36').Value = '  +0.48%  '

$ws.Range('D37').Value = '''2.764'
$ws.Range('E37').Value = '  +1.97%  '

$ws.Range('D38').Value = '''0.01652'
$ws.Range('E38').Value = '  +1.66%  '

$ws.Range('D39').Value = '1.109.29'
$ws.Range('E39').Value = '  -0.13%  '

$ws.Range('D40').Value = '''6.117'

$ws.Range('D41').Value = '''0.8825'
$ws.Range('E41').Value = '  +0.72%  '

$ws.Range('E42').Value = '  -0.18%  '

$ws.Range('D43').Value = '''101.39'
$ws.Range('E43').Value = '  +0.65%  '

$ws.Range('D44').Value = '1.842.39'

$ws.Range('E45').Value = '  -1.25%  '

$ws.Range('D46').Value = '''57.61'
$ws.Range('E46').Value = '  +1.13%  '

$ws.Range('D47').Value = '''8.125'
$ws.Range('E47').Value = '  -1.00%  '

$ws.Range('D48').Value = '''1.002'
$ws.Range('E48').Value = '  -0.47%  '

$ws.Range('E49').Value = '  +0.26%  '

$ws.Range('D50').Value = '''0.4300'
$ws.Range('E50').Value = '  +0.20%  '

$ws.Range('D51').Value = '''6.043'
$ws.Range('E51').Value = '  -0.74%  '
